# Weekly update: a new price record (row) for "Macroferia Regional de Talca -
# Arveja Verde" is inserted above the existing row 85, pushing the previous
# rows 85-87 down to 86-88.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 85 (shifts old rows 85-87 down to 86-88).
$ws.Rows("85:85").Insert()

# Populate the new row 85 with the latest weekly record.
$ws.Range("A85").Value = 5
$ws.Range("B85").Value = "Macroferia Regional de Talca"
$ws.Range("C85").Value = "Maule"
$ws.Range("D85").Value = 44595
$ws.Range("E85").Value = 7
$ws.Range("F85").Value = 100112022
$ws.Range("G85").Value = "Arveja Verde"
$ws.Range("H85").Value = "Perfection"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 250
$ws.Range("K85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = 30000
$ws.Range("N85").Value = "$/saco 25 kilos"
$ws.Range("O85").Value = "Carahue"
$ws.Range("P85").Value = 1200
$ws.Range("Q85").Value = 25
$ws.Range("R85").Value = "Hortaliza"
